# Auto-optimize exam scheduling: dynamically adjusts exams per slot (1-4)
# to guarantee all courses are scheduled within date range.
# This reshuffles the exam-slot assignments on the Section_A and Section_B
# timetable sheets.

$wb = $excel.ActiveWorkbook

# ---- Section_A ----
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "Free"
$wsA.Range("D2").Value = "DS303"

$wsA.Range("B3").Value = "DS303"
$wsA.Range("C3").Value = "DS302"
$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "DS303"
$wsA.Range("F3").Value = "Free"

$wsA.Range("C5").Value = "CS307"
$wsA.Range("F5").Value = "DS302"

$wsA.Range("B6").Value = "Free"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "DS303 (Tutorial)"

$wsA.Range("D7").Value = "Free"
$wsA.Range("F7").Value = "CS307"

$wsA.Range("F8").Value = "DS302 (Tutorial)"

# ---- Section_B ----
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "DS303"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "CS307"

$wsB.Range("C3").Value = "DS303"
$wsB.Range("F3").Value = "DS303"

$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "CS307"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "DS302"
$wsB.Range("F5").Value = "DS302"

$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "DS303 (Tutorial)"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "DS302 (Tutorial)"

$wsB.Range("B7").Value = "DS302"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "Free"
$wsB.Range("F7").Value = "Free"
